# Applies crypto price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "27.936.88"
Set-TextCell "E2" "  +1.49%  "
Set-TextCell "D3" "1.642.20"
Set-TextCell "E3" "  +0.78%  "
Set-TextCell "E4" "  +0.01%  "
Set-TextCell "D5" "212.99"
Set-TextCell "E5" "  +0.74%  "
Set-TextCell "D6" "0.524"
Set-TextCell "E6" "  +0.59%  "
Set-TextCell "E7" "  -0.02%  "
Set-TextCell "D8" "23.49"
Set-TextCell "E8" "  +2.43%  "
Set-TextCell "E9" "  -1.37%  "
Set-TextCell "E10" "  +0.63%  "
Set-TextCell "D11" "0.0882"
Set-TextCell "E11" "  +2.32%  "
Set-TextCell "D12" "1.874.81"
Set-TextCell "E12" "  +0.80%  "
Set-TextCell "D13" "1.642.48"
Set-TextCell "E14" "  +1.22%  "
Set-TextCell "E15" "  +2.61%  "
Set-TextCell "D16" "65.57"
Set-TextCell "E16" "  +0.88%  "
Set-TextCell "D17" "27.933.74"
Set-TextCell "E17" "  +1.56%  "
Set-TextCell "D18" "232.79"
Set-TextCell "E18" "  +1.81%  "
Set-TextCell "E19" "  +0.69%  "
Set-TextCell "E20" "  +0.95%  "
Set-TextCell "E21" "  +0.04%  "
Set-TextCell "D22" "10.50"
Set-TextCell "E22" "  -2.04%  "
Set-TextCell "D23" "4.37"
Set-TextCell "E23" "  +0.15%  "
Set-TextCell "E24" "  -1.78%  "
Set-TextCell "D25" "152.93"
Set-TextCell "E25" "  +2.68%  "
Set-TextCell "D26" "6.90"
Set-TextCell "E26" "  +0.59%  "
Set-TextCell "E27" "  +0.79%  "
Set-TextCell "E28" "  +0.08%  "
Set-TextCell "E29" "  -0.06%  "
Set-TextCell "E30" "  +0.94%  "
Set-TextCell "E31" "  +0.79%  "
Set-TextCell "E32" "  +2.88%  "
Set-TextCell "E33" "  +0.60%  "
Set-TextCell "D34" "1.409.64"
Set-TextCell "E34" "  -3.63%  "
Set-TextCell "E35" "  +2.68%  "
Set-TextCell "E36" "  +1.47%  "
Set-TextCell "E37" "  +1.75%  "
Set-TextCell "E38" "  +0.75%  "
Set-TextCell "D39" "0.562"
Set-TextCell "E39" "  +0.71%  "
Set-TextCell "D40" "0.924"
Set-TextCell "E40" "  +0.90%  "
Set-TextCell "E41" "  +1.07%  "
Set-TextCell "E42" "  -0.04%  "
Set-TextCell "D43" "67.43"
Set-TextCell "E43" "  -0.95%  "
Set-TextCell "E44" "  +6.84%  "
Set-TextCell "D45" "5.52"
Set-TextCell "E45" "  +3.00%  "
Set-TextCell "E46" "  +0.06%  "
Set-TextCell "D47" "1.783.76"
Set-TextCell "E47" "  +0.83%  "
Set-TextCell "D48" "87.79"
Set-TextCell "E48" "  +0.50%  "
Set-TextCell "E49" "  +0.84%  "
Set-TextCell "E50" "  +0.33%  "
Set-TextCell "D51" "7.62"
Set-TextCell "E51" "  -0.84%  "
